$d = $word.ActiveDocument

# Mapping of old "A×B=C" equation strings to new ones, applied in document order.
$r = $d.Content
$r.Find.Execute("97×49=4753", $true, $false, $false, $false, $false, $true, 1, $false, "74×88=6512", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("79×60=4740", $true, $false, $false, $false, $false, $true, 1, $false, "92×84=7728", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("69×70=4830", $true, $false, $false, $false, $false, $true, 1, $false, "47×30=1410", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("28×29=812", $true, $false, $false, $false, $false, $true, 1, $false, "82×85=6970", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("31×49=1519", $true, $false, $false, $false, $false, $true, 1, $false, "17×67=1139", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("21×93=1953", $true, $false, $false, $false, $false, $true, 1, $false, "23×74=1702", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("58×62=3596", $true, $false, $false, $false, $false, $true, 1, $false, "65×85=5525", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("51×84=4284", $true, $false, $false, $false, $false, $true, 1, $false, "53×27=1431", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("23×23=529", $true, $false, $false, $false, $false, $true, 1, $false, "69×83=5727", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("82×55=4510", $true, $false, $false, $false, $false, $true, 1, $false, "50×71=3550", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("43×66=2838", $true, $false, $false, $false, $false, $true, 1, $false, "68×96=6528", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("17×47=799", $true, $false, $false, $false, $false, $true, 1, $false, "99×83=8217", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("85×16=1360", $true, $false, $false, $false, $false, $true, 1, $false, "66×84=5544", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("25×55=1375", $true, $false, $false, $false, $false, $true, 1, $false, "25×30=750", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("11×51=561", $true, $false, $false, $false, $false, $true, 1, $false, "40×70=2800", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("38×20=760", $true, $false, $false, $false, $false, $true, 1, $false, "13×34=442", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("70×31=2170", $true, $false, $false, $false, $false, $true, 1, $false, "38×45=1710", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("27×90=2430", $true, $false, $false, $false, $false, $true, 1, $false, "95×74=7030", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("84×19=1596", $true, $false, $false, $false, $false, $true, 1, $false, "13×74=962", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("70×52=3640", $true, $false, $false, $false, $false, $true, 1, $false, "71×60=4260", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("14×25=350", $true, $false, $false, $false, $false, $true, 1, $false, "89×19=1691", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("73×30=2190", $true, $false, $false, $false, $false, $true, 1, $false, "31×83=2573", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("33×91=3003", $true, $false, $false, $false, $false, $true, 1, $false, "77×47=3619", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("19×83=1577", $true, $false, $false, $false, $false, $true, 1, $false, "76×75=5700", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("75×26=1950", $true, $false, $false, $false, $false, $true, 1, $false, "12×15=180", 2) | Out-Null
